$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new sequential values: P1=14, Q1=15
# matching the existing header style (bold, centered, thin border, top-aligned)
$ws.Range("P1").Value = 14
$ws.Range("P1").Font.Bold = $true
$ws.Range("P1").HorizontalAlignment = -4108
$ws.Range("P1").VerticalAlignment = -4160
$ws.Range("P1").Borders.LineStyle = 1

$ws.Range("Q1").Value = 15
$ws.Range("Q1").Font.Bold = $true
$ws.Range("Q1").HorizontalAlignment = -4108
$ws.Range("Q1").VerticalAlignment = -4160
$ws.Range("Q1").Borders.LineStyle = 1

# For data rows 2..25:
#  - swap values in columns I/K and M/O (1<->2)
#  - add new columns P and Q, both with value 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P -> 2
    $ws.Cells.Item($r, 17).Value = 2  # Q -> 2
}
